# Weekly update: insert a new observation as row 49, pushing the rest of the
# table (old rows 49-150) down by one row (new rows 50-151).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 49; existing rows 49:150 shift down to 50:151.
$ws.Rows("49:49").Insert()

# Populate the new row 49 with the new weekly data point.
$ws.Cells.Item(49, 1).Value  = 11
$ws.Cells.Item(49, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value  = "Bíobío"
$ws.Cells.Item(49, 4).Value  = 44519
$ws.Cells.Item(49, 5).Value  = 8
$ws.Cells.Item(49, 6).Value  = 100114013
$ws.Cells.Item(49, 7).Value  = "Zanahoria"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 250
$ws.Cells.Item(49, 11).Value = 7500
$ws.Cells.Item(49, 12).Value = 8000
$ws.Cells.Item(49, 13).Value = 7800
$ws.Cells.Item(49, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 390
$ws.Cells.Item(49, 17).Value = 20
$ws.Cells.Item(49, 18).Value = "Hortaliza"
